# Apply the PNAD 2009 "roubofurto" correction:
#  1. Rename the mislabeled header cell B2 from "unnamed: 1_level_1" to "total"
#  2. Remove the two label-only separator rows ("situação do domicílio" and
#     "grandes regiões e unidades da federação") which had no data of their own;
#     deleting them shifts the numeric data up so it lines up with the correct
#     region/label row (fixing the previous off-by-one misalignment) and shrinks
#     the shared-strings table accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the mislabeled header in B2 ("unnamed: 1_level_1" -> "total")
$ws.Range("B2").Value = "total"

# Row 5 is "situação do domicílio" (label only, no numeric data) - remove it.
# This shifts "urbana"/"rural" (and everything below) up by one row.
$ws.Rows.Item(5).Delete()

# After the previous deletion, "grandes regiões e unidades da federação"
# (label only, no numeric data) is now row 7 - remove it too.
# This shifts "norte" and all the states/regions below up by one more row.
$ws.Rows.Item(7).Delete()
